$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Sat Feb 17 22:55:08 EST 2024"
$ws.Range("B3").Value = "Sat Feb 17 22:55:20 EST 2024"
$ws.Range("B5").Value = "Sat Feb 17 22:55:31 EST 2024"
$ws.Range("B6").Value = "Sat Feb 17 22:55:43 EST 2024"
$ws.Range("B7").Value = "Sat Feb 17 22:55:55 EST 2024"
